# This workbook logs weekly wholesale mango prices. A new week of
# observations (date 44505 = 2021-11-05) is inserted at the top of the
# data block (rows 669-671), pushing all the subsequent rows down by
# three positions (old 669-773 become new 672-776).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 669, shifting
# everything from 669 downward (this also grows the used range / the
# sheet's <dimension> to A1:T776 automatically, matching the target).
$ws.Rows("669:671").Insert()

# Columns that are constant across the whole data block.
$commonA = 6
$commonB = "Mercado Mayorista Lo Valledor de Santiago"
$commonC = "Metropolitana"
$commonD = 44505
$commonE = 13
$commonF = "Fruta"
$commonG = 100108
$commonH = "Tropicales y subtropicales"
$commonI = 100108002
$commonJ = "Mango"
$commonK = "Sin especificar"
$commonQ = "`$/bandeja 4 kilos"
$commonT = 4

# Per-row values: Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Origen, Precio $/Kg.
$rowData = @{
    669 = @("Especial", 706,  6000, 7000, 6323, "Perú", 1581)
    670 = @("Primera",  1486, 5500, 7000, 6026, "Perú", 1506)
    671 = @("Segunda",  1426, 5000, 7000, 5585, "Perú", 1396)
}

foreach ($r in 669..671) {
    $vals = $rowData[$r]

    $ws.Cells.Item($r, 1).Value2  = $commonA
    $ws.Cells.Item($r, 2).Value2  = $commonB
    $ws.Cells.Item($r, 3).Value2  = $commonC
    $ws.Cells.Item($r, 4).Value2  = $commonD
    $ws.Cells.Item($r, 5).Value2  = $commonE
    $ws.Cells.Item($r, 6).Value2  = $commonF
    $ws.Cells.Item($r, 7).Value2  = $commonG
    $ws.Cells.Item($r, 8).Value2  = $commonH
    $ws.Cells.Item($r, 9).Value2  = $commonI
    $ws.Cells.Item($r, 10).Value2 = $commonJ
    $ws.Cells.Item($r, 11).Value2 = $commonK
    $ws.Cells.Item($r, 12).Value2 = $vals[0]
    $ws.Cells.Item($r, 13).Value2 = $vals[1]
    $ws.Cells.Item($r, 14).Value2 = $vals[2]
    $ws.Cells.Item($r, 15).Value2 = $vals[3]
    $ws.Cells.Item($r, 16).Value2 = $vals[4]
    $ws.Cells.Item($r, 17).Value2 = $commonQ
    $ws.Cells.Item($r, 18).Value2 = $vals[5]
    $ws.Cells.Item($r, 19).Value2 = $vals[6]
    $ws.Cells.Item($r, 20).Value2 = $commonT
}
